$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = [double]"24.13000000000033"
$ws.Range("H2").Value = [double]"0.007445307440328675"
$ws.Range("I2").Value = [double]"0.007445307440328675"
$ws.Range("L2").Value = [double]"34.84543851657722"
$ws.Range("M2").Value = "[9.559682795670042, 60.131194237484394]"
$ws.Range("N2").Value = [double]"0.007999541567619062"
$ws.Range("O2").Value = [double]"0.007999541567619062"
$ws.Range("P2").Value = [double]"1.62897396852804"
$ws.Range("Q2").Value = "[0.672973801669885, 2.5849741353861955]"
$ws.Range("R2").Value = [double]"0.001295700348463269"
$ws.Range("S2").Value = [double]"0.001295700348463269"
$ws.Range("T2").Value = [double]"60.25361881599026"
$ws.Range("U2").Value = "[45.22273939370312, 75.2844982382774]"
$ws.Range("V2").Value = [double]"2.692059908326883e-10"
$ws.Range("W2").Value = [double]"2.692059908326883e-10"
$ws.Range("X2").Value = [double]"17.87407407407432"
$ws.Range("Y2").Value = [double]"14.20264264264284"
$ws.Range("Z2").Value = [double]"21.5455055055058"

$ws.Range("B3").Value = [double]"0"
$ws.Range("F3").Value = [double]"24.13000000000033"
$ws.Range("H3").Value = [double]"0.01822275359122094"
$ws.Range("I3").Value = [double]"0.01822275359122094"
$ws.Range("L3").Value = [double]"36.93357835851143"
$ws.Range("M3").Value = "[4.2906647800868996, 69.57649193693595]"
$ws.Range("N3").Value = [double]"0.02747210930424115"
$ws.Range("O3").Value = [double]"0.02747210930424115"
$ws.Range("P3").Value = [double]"1.125815985971117"
$ws.Range("Q3").Value = "[-0.14465791998511612, 2.396289891927349]"
$ws.Range("R3").Value = [double]"0.08103996967428384"
$ws.Range("S3").Value = [double]"0.08103996967428384"
$ws.Range("T3").Value = [double]"71.26833453916416"
$ws.Range("U3").Value = "[53.49074293133708, 89.04592614699123]"
$ws.Range("V3").Value = [double]"2.687887690200341e-10"
$ws.Range("W3").Value = [double]"2.687887690200341e-10"
$ws.Range("X3").Value = [double]"19.80640640640668"
$ws.Range("Y3").Value = [double]"14.92726726726747"
$ws.Range("Z3").Value = [double]"24.68554554554589"

$ws.Range("F4").Value = [double]"24.13000000000033"
$ws.Range("H4").Value = [double]"0.0004261584349132264"
$ws.Range("I4").Value = [double]"0.0004261584349132264"
$ws.Range("L4").Value = [double]"47.275793750324"
$ws.Range("M4").Value = "[18.836867941354456, 75.71471955929354]"
$ws.Range("N4").Value = [double]"0.001651796124107818"
$ws.Range("O4").Value = [double]"0.001651796124107818"
$ws.Range("P4").Value = [double]"1.30192127986604"
$ws.Range("Q4").Value = "[0.6352369529781159, 1.9686056067539637]"
$ws.Range("R4").Value = [double]"0.000286967815912087"
$ws.Range("S4").Value = [double]"0.000286967815912087"
$ws.Range("T4").Value = [double]"54.69232870975982"
$ws.Range("U4").Value = "[39.0222986924448, 70.36235872707485]"
$ws.Range("V4").Value = [double]"9.206604145717279e-09"
$ws.Range("W4").Value = [double]"9.206604145717279e-09"
$ws.Range("X4").Value = [double]"19.13009009009036"
$ws.Range("Y4").Value = [double]"16.56974974974998"
$ws.Range("Z4").Value = [double]"21.69043043043073"

$ws.Range("F5").Value = [double]"24.13000000000033"
$ws.Range("H5").Value = [double]"0.0001393026006277776"
$ws.Range("I5").Value = [double]"0.0001393026006277776"
$ws.Range("L5").Value = [double]"51.25566456310855"
$ws.Range("M5").Value = "[21.54779721026395, 80.96353191595314]"
$ws.Range("N5").Value = [double]"0.001142430772846037"
$ws.Range("O5").Value = [double]"0.001142430772846037"
$ws.Range("P5").Value = [double]"1.025184389459732"
$ws.Range("Q5").Value = "[0.42139481039142357, 1.628973968528041]"
$ws.Range("R5").Value = [double]"0.001342371809438303"
$ws.Range("S5").Value = [double]"0.001342371809438303"
$ws.Range("T5").Value = [double]"62.33028836894802"
$ws.Range("U5").Value = "[46.67293234624644, 77.9876443916496]"
$ws.Range("V5").Value = [double]"3.246065638506934e-10"
$ws.Range("W5").Value = [double]"3.246065638506934e-10"
$ws.Range("X5").Value = [double]"20.19287287287315"
$ws.Range("Y5").Value = [double]"17.87407407407432"
$ws.Range("Z5").Value = [double]"22.51167167167198"

$ws.Range("F6").Value = [double]"24.13000000000033"
$ws.Range("H6").Value = [double]"0.0002108185848478383"
$ws.Range("I6").Value = [double]"0.0002108185848478383"
$ws.Range("L6").Value = [double]"50.14808706060595"
$ws.Range("M6").Value = "[22.682787822330212, 77.61338629888168]"
$ws.Range("N6").Value = [double]"0.000625997030035963"
$ws.Range("O6").Value = [double]"0.000625997030035963"
$ws.Range("P6").Value = [double]"1.377394977249579"
$ws.Range("Q6").Value = "[0.7107106503616549, 2.0440793041375027]"
$ws.Range("R6").Value = [double]"0.0001407307904610633"
$ws.Range("S6").Value = [double]"0.0001407307904610633"
$ws.Range("T6").Value = [double]"64.15771277720145"
$ws.Range("U6").Value = "[48.37814290570071, 79.9372826487022]"
$ws.Range("V6").Value = [double]"1.832121121481123e-10"
$ws.Range("W6").Value = [double]"1.832121121481123e-10"
$ws.Range("X6").Value = [double]"18.8402402402405"
$ws.Range("Y6").Value = [double]"16.27989989990013"
$ws.Range("Z6").Value = [double]"21.40058058058088"

$ws.Range("F7").Value = [double]"24.13000000000033"
$ws.Range("H7").Value = [double]"0.003013510730262126"
$ws.Range("I7").Value = [double]"0.003013510730262126"
$ws.Range("L7").Value = [double]"38.45779374406717"
$ws.Range("M7").Value = "[12.474918640172234, 64.4406688479621]"
$ws.Range("N7").Value = [double]"0.004621264171096096"
$ws.Range("O7").Value = [double]"0.004621264171096096"
$ws.Range("P7").Value = [double]"1.452868674633118"
$ws.Range("Q7").Value = "[0.5849211547224247, 2.320816194543811]"
$ws.Range("R7").Value = [double]"0.001544503290096699"
$ws.Range("S7").Value = [double]"0.001544503290096699"
$ws.Range("T7").Value = [double]"37.60303634370829"
$ws.Range("U7").Value = "[22.53080670303111, 52.67526598438547]"
$ws.Range("V7").Value = [double]"8.466853261612783e-06"
$ws.Range("W7").Value = [double]"8.466853261612783e-06"
$ws.Range("X7").Value = [double]"18.55039039039065"
$ws.Range("Y7").Value = [double]"15.21711711711733"
$ws.Range("Z7").Value = [double]"21.88366366366396"

$ws.Range("F8").Value = [double]"23.46000000000023"
$ws.Range("H8").Value = [double]"1.110844504081943e-05"
$ws.Range("I8").Value = [double]"1.110844504081943e-05"
$ws.Range("L8").Value = [double]"50.45669041524739"
$ws.Range("M8").Value = "[25.09127687888421, 75.82210395161056]"
$ws.Range("N8").Value = [double]"0.0002286411538241495"
$ws.Range("O8").Value = [double]"0.0002286411538241495"
$ws.Range("P8").Value = [double]"0.6603948521059628"
$ws.Range("Q8").Value = "[0.16981581911296217, 1.1509738850989635]"
$ws.Range("R8").Value = [double]"0.009455502149745332"
$ws.Range("S8").Value = [double]"0.009455502149745332"
$ws.Range("T8").Value = [double]"53.95821907594169"
$ws.Range("U8").Value = "[40.74478903013253, 67.17164912175085]"
$ws.Range("V8").Value = [double]"1.626700996126829e-10"
$ws.Range("W8").Value = [double]"1.626700996126829e-10"
$ws.Range("X8").Value = [double]"20.99423423423444"
$ws.Range("Y8").Value = [double]"19.16252252252271"
$ws.Range("Z8").Value = [double]"22.82594594594617"

$ws.Range("B9").Value = [double]"1"
$ws.Range("F9").Value = [double]"23.46000000000023"
$ws.Range("H9").Value = [double]"2.58302617583972e-06"
$ws.Range("I9").Value = [double]"2.58302617583972e-06"
$ws.Range("L9").Value = [double]"53.93515461421167"
$ws.Range("M9").Value = "[28.616033564430168, 79.25427566399317]"
$ws.Range("N9").Value = [double]"9.337069441395229e-05"
$ws.Range("O9").Value = [double]"9.337069441395229e-05"
$ws.Range("P9").Value = [double]"0.6226580034141938"
$ws.Range("Q9").Value = "[0.15723686954904004, 1.0880791372793475]"
$ws.Range("R9").Value = [double]"0.009873137367608198"
$ws.Range("S9").Value = [double]"0.009873137367608198"
$ws.Range("T9").Value = [double]"55.41138821466592"
$ws.Range("U9").Value = "[42.34446219547229, 68.47831423385955]"
$ws.Range("V9").Value = [double]"5.701017435910671e-11"
$ws.Range("W9").Value = [double]"5.701017435910671e-11"
$ws.Range("X9").Value = [double]"21.13513513513534"
$ws.Range("Y9").Value = [double]"19.39735735735755"
$ws.Range("Z9").Value = [double]"22.87291291291313"

$ws.Range("B10").Value = [double]"0"
$ws.Range("F10").Value = [double]"23.46000000000023"
$ws.Range("H10").Value = [double]"0.1087888764827238"
$ws.Range("I10").Value = [double]"0.1087888764827238"
$ws.Range("L10").Value = [double]"29.91857618402186"
$ws.Range("M10").Value = "[-7.3639857159477025, 67.20113808399142]"
$ws.Range("N10").Value = [double]"0.1130235844132124"
$ws.Range("O10").Value = [double]"0.1130235844132124"
$ws.Range("P10").Value = [double]"0.798763297309117"
$ws.Range("Q10").Value = "[-2.32081619454381, 3.918342789162044]"
$ws.Range("R10").Value = [double]"0.608583029557269"
$ws.Range("S10").Value = [double]"0.608583029557269"
$ws.Range("T10").Value = [double]"62.73956433834068"
$ws.Range("U10").Value = "[43.011465678298464, 82.4676629983829]"
$ws.Range("V10").Value = [double]"7.798467027697598e-08"
$ws.Range("W10").Value = [double]"7.798467027697598e-08"
$ws.Range("X10").Value = [double]"20.4775975975978"
$ws.Range("Y10").Value = [double]"8.829789789789878"
$ws.Range("Z10").Value = [double]"32.12540540540571"

$ws.Range("F11").Value = [double]"23.46000000000023"
$ws.Range("H11").Value = [double]"0.0004191859398923192"
$ws.Range("I11").Value = [double]"0.0004191859398923192"
$ws.Range("L11").Value = [double]"48.12999955436025"
$ws.Range("M11").Value = "[17.69244925678217, 78.56754985193832]"
$ws.Range("N11").Value = [double]"0.002630018641599507"
$ws.Range("O11").Value = [double]"0.002630018641599507"
$ws.Range("P11").Value = [double]"1.062921238151502"
$ws.Range("Q11").Value = "[0.4339737599553466, 1.691868716347658]"
$ws.Range("R11").Value = [double]"0.00140607517553315"
$ws.Range("S11").Value = [double]"0.00140607517553315"
$ws.Range("T11").Value = [double]"63.12679237043573"
$ws.Range("U11").Value = "[47.31973394138508, 78.93385079948638]"
$ws.Range("V11").Value = [double]"2.979876345676757e-10"
$ws.Range("W11").Value = [double]"2.979876345676757e-10"
$ws.Range("X11").Value = [double]"19.49129129129148"
$ws.Range("Y11").Value = [double]"17.1429429429431"
$ws.Range("Z11").Value = [double]"21.83963963963985"

$ws.Range("F12").Value = [double]"23.46000000000023"
$ws.Range("H12").Value = [double]"8.796662047672044e-06"
$ws.Range("I12").Value = [double]"8.796662047672044e-06"
$ws.Range("L12").Value = [double]"53.99765674967724"
$ws.Range("M12").Value = "[27.90395983318595, 80.09135366616853]"
$ws.Range("N12").Value = [double]"0.0001377746041228978"
$ws.Range("O12").Value = [double]"0.0001377746041228978"
$ws.Range("P12").Value = [double]"1.088079137279347"
$ws.Range("Q12").Value = "[0.5723422051585008, 1.6038160694001942]"
$ws.Range("R12").Value = [double]"0.0001064666480921073"
$ws.Range("S12").Value = [double]"0.0001064666480921073"
$ws.Range("T12").Value = [double]"56.96799453066256"
$ws.Range("U12").Value = "[43.12961551430105, 70.80637354702407]"
$ws.Range("V12").Value = [double]"1.303228636118092e-10"
$ws.Range("W12").Value = [double]"1.303228636118092e-10"
$ws.Range("X12").Value = [double]"19.39735735735755"
$ws.Range("Y12").Value = [double]"17.47171171171189"
$ws.Range("Z12").Value = [double]"21.32300300300321"
